$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for fab95b9e row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-30 09:02:42"

# zh-cn sheet: Correspond Handoff Datetime (H4) and Correspond Handback DateTime (K4) for fab95b9e row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-30 09:02:37"
$wsZhCn.Range("K4").Value = "2016-08-30 09:02:55"

# de-de sheet: Correspond Handback DateTime (K4) for fab95b9e row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-30 09:03:12"
